$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"   = -6.939899999999994
    "C12"  = -11.0418
    "D14"  = -7.387200000000004
    "D26"  = -8.347399999999997
    "C27"  = -12.7943
    "D31"  = -8.160500000000001
    "C32"  = -13.3625
    "D35"  = -8.407299999999994
    "C36"  = -12.70210000000001
    "D37"  = -7.873999999999994
    "C38"  = -12.89399999999999
    "D45"  = -7.6271
    "C46"  = -14.64799999999999
    "D52"  = -7.621599999999998
    "C54"  = -13.03650000000001
    "C55"  = -13.6077
    "C56"  = -12.81909999999999
    "D57"  = -8.364299999999998
    "C67"  = -10.68750000000001
    "C69"  = -11.87479999999999
    "C72"  = -11.6041
    "D81"  = -6.965799999999994
    "C83"  = -13.56920000000001
    "D83"  = -9.0039
    "C86"  = -13.45339999999999
    "C91"  = -10.3491
    "C93"  = -10.83910000000001
    "C99"  = -13.01299999999999
    "D100" = -7.973600000000004
    "D102" = -7.859299999999996
}

foreach ($cellAddr in $updates.Keys) {
    $ws.Range($cellAddr).Value = $updates[$cellAddr]
}
